# Updates cryptos list values (price/volume columns) and re-orders a few
# coin rows (B/C columns), matching the "Updated cryptos list" data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.993.18"
$ws.Range("E2").Value = "  +4.56%  "
$ws.Range("D3").Value = "3.443.03"
$ws.Range("E3").Value = "  +3.76%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.56"
$ws.Range("E5").Value = "  +5.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.47"
$ws.Range("E6").Value = "  +7.20%  "
$ws.Range("E7").Value = "  +2.43%  "
$ws.Range("D8").Value = "3.437.89"
$ws.Range("E8").Value = "  +3.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +2.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.646"
$ws.Range("E11").Value = "  +2.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.32"
$ws.Range("E12").Value = "  +5.52%  "
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.44"
$ws.Range("E14").Value = "  +4.67%  "
$ws.Range("D15").Value = "3.995.13"
$ws.Range("E15").Value = "  +3.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.67"
$ws.Range("E16").Value = "  +3.38%  "
$ws.Range("D17").Value = "3.455.54"
$ws.Range("E17").Value = "  +4.18%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.120"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "66.918.82"
$ws.Range("E19").Value = "  +4.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.07"
$ws.Range("E20").Value = "  +3.57%  "
$ws.Range("E21").Value = "  +3.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "482.67"
$ws.Range("E22").Value = "  +6.92%  "
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.30"
$ws.Range("E23").Value = "  +5.77%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.85"
$ws.Range("E24").Value = "  +22.89%  "
$ws.Range("E25").Value = "  +8.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.76"
$ws.Range("E26").Value = "  +3.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.03"
$ws.Range("E27").Value = "  +3.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.94"
$ws.Range("E28").Value = "  +3.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.16"
$ws.Range("E29").Value = "  +7.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.32"
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.13"
$ws.Range("E31").Value = "  +9.37%  "
$ws.Range("E32").Value = "  +2.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "64.10"
$ws.Range("E33").Value = "  +5.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "590.17"
$ws.Range("E34").Value = "  +4.70%  "
$ws.Range("E35").Value = "  +4.87%  "
$ws.Range("E36").Value = "  +5.98%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.57"
$ws.Range("E38").Value = "  +1.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.44"
$ws.Range("E39").Value = "  +3.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.385"
$ws.Range("E40").Value = "  +5.27%  "
$ws.Range("D41").Value = "0.0₃0770"
$ws.Range("E41").Value = "  +5.95%  "
$ws.Range("D42").Value = "3.193.70"
$ws.Range("E42").Value = "  +4.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.91"
$ws.Range("E43").Value = "  +5.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0430"
$ws.Range("E44").Value = "  +4.28%  "
$ws.Range("E45").Value = "  +4.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.77"
$ws.Range("E46").Value = "  +21.40%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.22"
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.135"
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.71"
$ws.Range("E50").Value = "  +7.36%  "
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.19"
$ws.Range("E51").Value = "  +9.46%  "
